$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-97 (96 data rows) of the Actual Production Wind sheet are being
# "retrained" - the timestamp in column A moves forward 4 days (from
# 2025-09-16 to 2025-09-20 in Excel serial date terms) and the wind
# production values in column B are replaced with the newly observed data.
$data = New-Object 'object[,]' 96,2
$data[0,0] = 45923.01041666666
$data[0,1] = 652
$data[1,0] = 45923.02083333334
$data[1,1] = 620
$data[2,0] = 45923.03125
$data[2,1] = 574
$data[3,0] = 45923.04166666666
$data[3,1] = 553
$data[4,0] = 45923.05208333334
$data[4,1] = 523
$data[5,0] = 45923.0625
$data[5,1] = 491
$data[6,0] = 45923.07291666666
$data[6,1] = 468
$data[7,0] = 45923.08333333334
$data[7,1] = 442
$data[8,0] = 45923.09375
$data[8,1] = 417
$data[9,0] = 45923.10416666666
$data[9,1] = 410
$data[10,0] = 45923.11458333334
$data[10,1] = 388
$data[11,0] = 45923.125
$data[11,1] = 349
$data[12,0] = 45923.13541666666
$data[12,1] = 300
$data[13,0] = 45923.14583333334
$data[13,1] = 273
$data[14,0] = 45923.15625
$data[14,1] = 242
$data[15,0] = 45923.16666666666
$data[15,1] = 223
$data[16,0] = 45923.17708333334
$data[16,1] = 196
$data[17,0] = 45923.1875
$data[17,1] = 179
$data[18,0] = 45923.19791666666
$data[18,1] = 177
$data[19,0] = 45923.20833333334
$data[19,1] = 162
$data[20,0] = 45923.21875
$data[20,1] = 172
$data[21,0] = 45923.22916666666
$data[21,1] = 175
$data[22,0] = 45923.23958333334
$data[22,1] = 169
$data[23,0] = 45923.25
$data[23,1] = 166
$data[24,0] = 45923.26041666666
$data[24,1] = 172
$data[25,0] = 45923.27083333334
$data[25,1] = 182
$data[26,0] = 45923.28125
$data[26,1] = 181
$data[27,0] = 45923.29166666666
$data[27,1] = 0
$data[28,0] = 45923.30208333334
$data[28,1] = 0
$data[29,0] = 45923.3125
$data[29,1] = 0
$data[30,0] = 45923.32291666666
$data[30,1] = 0
$data[31,0] = 45923.33333333334
$data[31,1] = 0
$data[32,0] = 45923.34375
$data[32,1] = 0
$data[33,0] = 45923.35416666666
$data[33,1] = 0
$data[34,0] = 45923.36458333334
$data[34,1] = 0
$data[35,0] = 45923.375
$data[35,1] = 0
$data[36,0] = 45923.38541666666
$data[36,1] = 0
$data[37,0] = 45923.39583333334
$data[37,1] = 0
$data[38,0] = 45923.40625
$data[38,1] = 0
$data[39,0] = 45923.41666666666
$data[39,1] = 0
$data[40,0] = 45923.42708333334
$data[40,1] = 0
$data[41,0] = 45923.4375
$data[41,1] = 0
$data[42,0] = 45923.44791666666
$data[42,1] = 0
$data[43,0] = 45923.45833333334
$data[43,1] = 0
$data[44,0] = 45923.46875
$data[44,1] = 0
$data[45,0] = 45923.47916666666
$data[45,1] = 0
$data[46,0] = 45923.48958333334
$data[46,1] = 0
$data[47,0] = 45923.5
$data[47,1] = 0
$data[48,0] = 45923.51041666666
$data[48,1] = 0
$data[49,0] = 45923.52083333334
$data[49,1] = 0
$data[50,0] = 45923.53125
$data[50,1] = 0
$data[51,0] = 45923.54166666666
$data[51,1] = 0
$data[52,0] = 45923.55208333334
$data[52,1] = 0
$data[53,0] = 45923.5625
$data[53,1] = 0
$data[54,0] = 45923.57291666666
$data[54,1] = 0
$data[55,0] = 45923.58333333334
$data[55,1] = 0
$data[56,0] = 45923.59375
$data[56,1] = 0
$data[57,0] = 45923.60416666666
$data[57,1] = 0
$data[58,0] = 45923.61458333334
$data[58,1] = 0
$data[59,0] = 45923.625
$data[59,1] = 0
$data[60,0] = 45923.63541666666
$data[60,1] = 0
$data[61,0] = 45923.64583333334
$data[61,1] = 0
$data[62,0] = 45923.65625
$data[62,1] = 0
$data[63,0] = 45923.66666666666
$data[63,1] = 0
$data[64,0] = 45923.67708333334
$data[64,1] = 0
$data[65,0] = 45923.6875
$data[65,1] = 0
$data[66,0] = 45923.69791666666
$data[66,1] = 0
$data[67,0] = 45923.70833333334
$data[67,1] = 0
$data[68,0] = 45923.71875
$data[68,1] = 0
$data[69,0] = 45923.72916666666
$data[69,1] = 0
$data[70,0] = 45923.73958333334
$data[70,1] = 0
$data[71,0] = 45923.75
$data[71,1] = 0
$data[72,0] = 45923.76041666666
$data[72,1] = 0
$data[73,0] = 45923.77083333334
$data[73,1] = 0
$data[74,0] = 45923.78125
$data[74,1] = 0
$data[75,0] = 45923.79166666666
$data[75,1] = 0
$data[76,0] = 45923.80208333334
$data[76,1] = 0
$data[77,0] = 45923.8125
$data[77,1] = 0
$data[78,0] = 45923.82291666666
$data[78,1] = 0
$data[79,0] = 45923.83333333334
$data[79,1] = 0
$data[80,0] = 45923.84375
$data[80,1] = 0
$data[81,0] = 45923.85416666666
$data[81,1] = 0
$data[82,0] = 45923.86458333334
$data[82,1] = 0
$data[83,0] = 45923.875
$data[83,1] = 0
$data[84,0] = 45923.88541666666
$data[84,1] = 0
$data[85,0] = 45923.89583333334
$data[85,1] = 0
$data[86,0] = 45923.90625
$data[86,1] = 0
$data[87,0] = 45923.91666666666
$data[87,1] = 0
$data[88,0] = 45923.92708333334
$data[88,1] = 0
$data[89,0] = 45923.9375
$data[89,1] = 0
$data[90,0] = 45923.94791666666
$data[90,1] = 0
$data[91,0] = 45923.95833333334
$data[91,1] = 0
$data[92,0] = 45923.96875
$data[92,1] = 0
$data[93,0] = 45923.97916666666
$data[93,1] = 0
$data[94,0] = 45923.98958333334
$data[94,1] = 0
$data[95,0] = 45924
$data[95,1] = 0

$ws.Range("A2:B97").Value = $data

Write-Output "Updated A2:B97 with new timestamps and production values"
